# Update the "dSF" column (F) values for several rows to reflect the
# re-pulled / recalculated data, per commit message:
# "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -3
$ws.Range("F5").Value = -4
$ws.Range("F6").Value = -2
$ws.Range("F7").Value = -6
$ws.Range("F9").Value = -4
$ws.Range("F11").Value = 0
$ws.Range("F14").Value = -7
$ws.Range("F15").Value = -7
$ws.Range("F16").Value = 1
$ws.Range("F21").Value = 0
$ws.Range("F25").Value = 0
